# Apply updated crypto price/volume data (GitHub Actions refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "25.913.46"
$ws.Range("E2").Value = "  -1.37%  "
$ws.Range("D3").Value = "1.637.15"
$ws.Range("E3").Value = "  -0.66%  "
$ws.Range("E4").Value = "  +0.15%  "
$ws.Range("D5").Value = "'215.48"
$ws.Range("E6").Value = "  +0.33%  "
$ws.Range("E7").Value = "  +0.12%  "
$ws.Range("E8").Value = "  -0.95%  "
$ws.Range("E9").Value = "  -0.16%  "
$ws.Range("E10").Value = "  -1.86%  "
$ws.Range("E11").Value = "  -0.07%  "
$ws.Range("E12").Value = "  -0.58%  "
$ws.Range("D13").Value = "1.863.42"
$ws.Range("E13").Value = "  -0.65%  "
$ws.Range("D14").Value = "1.627.43"
$ws.Range("E14").Value = "  -0.49%  "
$ws.Range("E15").Value = "  -0.75%  "
$ws.Range("D16").Value = "0.0₃0763"
$ws.Range("E16").Value = "  -0.42%  "
$ws.Range("D17").Value = "'62.85"
$ws.Range("E17").Value = "  -0.85%  "
$ws.Range("D18").Value = "25.954.69"
$ws.Range("E18").Value = "  -1.24%  "
$ws.Range("E19").Value = "  +0.15%  "
$ws.Range("D20").Value = "'193.17"
$ws.Range("E20").Value = "  -1.20%  "
$ws.Range("E21").Value = "  -2.18%  "
$ws.Range("E22").Value = "  -1.73%  "
$ws.Range("E23").Value = "  -0.87%  "
$ws.Range("E24").Value = "  +4.74%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("D26").Value = "'143.27"
$ws.Range("E27").Value = "  +0.08%  "
$ws.Range("D28").Value = "'6.89"
$ws.Range("E28").Value = "  -0.92%  "
$ws.Range("E29").Value = "  -0.64%  "
$ws.Range("E30").Value = "  -0.72%  "
$ws.Range("E31").Value = "  -0.74%  "
$ws.Range("E32").Value = "  -2.22%  "
$ws.Range("D33").Value = "'3.25"
$ws.Range("E33").Value = "  -0.31%  "
$ws.Range("D34").Value = "'1.53"
$ws.Range("E34").Value = "  -4.34%  "
$ws.Range("E36").Value = "  -1.50%  "
$ws.Range("D37").Value = "1.133.59"
$ws.Range("E37").Value = "  -0.44%  "
$ws.Range("D38").Value = "'0.544"
$ws.Range("E38").Value = "  -1.81%  "
$ws.Range("E39").Value = "  -1.61%  "
$ws.Range("E40").Value = "  -0.72%  "
$ws.Range("E41").Value = "  -0.84%  "
$ws.Range("D42").Value = "'99.32"
$ws.Range("E42").Value = "  -1.27%  "
$ws.Range("D43").Value = "'0.796"
$ws.Range("E43").Value = "  -0.56%  "
$ws.Range("D44").Value = "1.773.02"
$ws.Range("E44").Value = "  -0.61%  "
$ws.Range("E45").Value = "  +2.80%  "
$ws.Range("D46").Value = "'56.57"
$ws.Range("E46").Value = "  -1.13%  "
$ws.Range("D47").Value = "'0.0529"
$ws.Range("E47").Value = "  +2.21%  "
$ws.Range("E48").Value = "  -0.07%  "
$ws.Range("D49").Value = "'7.67"
$ws.Range("E49").Value = "  -0.57%  "
$ws.Range("E50").Value = "  -0.91%  "
$ws.Range("D51").Value = "'0.0958"
$ws.Range("E51").Value = "  -1.45%  "
